# Fruta / hortaliza, semanal
# Insert a new weekly record at row 83 (pushing existing rows 83:155 down to 84:156)
# and populate it with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(83).Insert()

$ws.Range("A83").Value = 10
$ws.Range("B83").Value = "Vega Modelo de Temuco"
$ws.Range("C83").Value = "La Araucanía"
$ws.Range("D83").Value = 45236
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = 100112010
$ws.Range("G83").Value = "Achicoria"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 100
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = 10000
$ws.Range("N83").Value = "`$/caja 18 unidades"
$ws.Range("O83").Value = "Región de O'Higgins"
$ws.Range("P83").Value = 556
$ws.Range("Q83").Value = 18
$ws.Range("R83").Value = "Hortaliza"
